$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.235.67"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "3.074.20"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'521.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").Value = "'135.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.63%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.073.54"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").Value = "'0.465"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.36%  "
$ws.Range("D10").Value = "'7.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("E12").Value = "  +1.87%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "3.604.36"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "'25.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").Value = "'0.0000161"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "57.293.31"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "3.072.46"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").Value = "'5.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D20").Value = "'12.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("D21").Value = "'7.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").Value = "'350.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.25%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'69.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").Value = "'0.498"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.17%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.165"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "0.0₃0864"
$ws.Range("E28").Value = "  -6.28%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "'7.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("D32").Value = "'5.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.72%  "
$ws.Range("D33").Value = "'20.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "'4.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.87%  "
$ws.Range("D35").Value = "'158.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  -4.98%  "
$ws.Range("D37").Value = "'6.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("D38").Value = "'25.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("D39").Value = "'1.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("D40").Value = "'0.0656"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("E41").Value = "  -5.34%  "
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("D43").Value = "'0.692"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "2.405.58"
$ws.Range("E44").Value = "  +5.61%  "
$ws.Range("D45").Value = "'36.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "3.113.07"
$ws.Range("E47").Value = "  -1.50%  "
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'0.945"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.41%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'5.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").Value = "'19.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.18%  "
